$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.767.38'
$ws.Range("E2").Value = '  -1.34%  '

$ws.Range("D3").Value = '1.595.77'
$ws.Range("E3").Value = '  -2.19%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.61'
$ws.Range("E5").Value = '  -2.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.507'
$ws.Range("E6").Value = '  -2.02%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("E8").Value = '  -1.88%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0619'
$ws.Range("E9").Value = '  -0.66%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.61'
$ws.Range("E10").Value = '  -2.44%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0833'
$ws.Range("E11").Value = '  -1.90%  '

$ws.Range("D12").Value = '1.818.38'
$ws.Range("E12").Value = '  -2.23%  '

$ws.Range("D13").Value = '1.583.69'
$ws.Range("E13").Value = '  -2.80%  '

$ws.Range("E14").Value = '  -1.25%  '

$ws.Range("E15").Value = '  -2.04%  '

$ws.Range("D16").Value = '26.770.23'
$ws.Range("E16").Value = '  -1.26%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.53'
$ws.Range("E17").Value = '  -3.24%  '

$ws.Range("D18").Value = '0.0₃0730'
$ws.Range("E18").Value = '  -0.30%  '

$ws.Range("E19").Value = '  -2.31%  '

$ws.Range("E20").Value = '  -0.03%  '

$ws.Range("E21").Value = '  -1.66%  '

$ws.Range("E22").Value = '  -2.55%  '

$ws.Range("E23").Value = '  -6.57%  '

$ws.Range("E24").Value = '  -2.73%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.59'
$ws.Range("E25").Value = '  -0.42%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.48'
$ws.Range("E26").Value = '  +1.55%  '

$ws.Range("E27").Value = '  -0.01%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.113'
$ws.Range("E28").Value = '  -4.59%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.32'
$ws.Range("E29").Value = '  -1.74%  '

$ws.Range("E30").Value = '  -0.55%  '

$ws.Range("E31").Value = '  -2.72%  '

$ws.Range("E32").Value = '  -2.73%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.674'
$ws.Range("E33").Value = '  +24.26%  '

$ws.Range("E34").Value = '  -2.10%  '

$ws.Range("D35").Value = '1.311.44'
$ws.Range("E35").Value = '  +0.36%  '

$ws.Range("E36").Value = '  -2.97%  '

$ws.Range("E37").Value = '  -0.79%  '

$ws.Range("E38").Value = '  -1.09%  '

$ws.Range("E39").Value = '  -2.87%  '

$ws.Range("E40").Value = '  -0.04%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.789'
$ws.Range("E41").Value = '  -1.83%  '

$ws.Range("E42").Value = '  -4.29%  '

$ws.Range("E43").Value = '  +0.26%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.79'
$ws.Range("E44").Value = '  +1.05%  '

$ws.Range("D45").Value = '1.731.80'
$ws.Range("E45").Value = '  -2.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.01'
$ws.Range("E46").Value = '  -1.77%  '

$ws.Range("E47").Value = '  +0.74%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.810'
$ws.Range("E48").Value = '  -0.38%  '

$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0102'
$ws.Range("E49").Value = '  -4.42%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0510'
$ws.Range("E50").Value = '  -0.75%  '

$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0978'
$ws.Range("E51").Value = '  +2.89%  '
